$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25: "What is the mean of 5/ 11?" -> " What is the mean of the set {0, 5, 10}?" ---
$ws.Range("A25").Value = " What is the mean of the set {0, 5, 10}?"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 2.5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = "C"
$ws.Range("G25:J25").ClearContents()

# --- Row 26: "What is the median of 3/ 5?" -> "What is the mean of the set {3, 3, 3, 3, 3, 0, 6}?" ---
$ws.Range("A26").Value = "What is the mean of the set {3, 3, 3, 3, 3, 0, 6}?"
$ws.Range("B26").Value = 1.5
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 6
$ws.Range("F26").Value = "A"
$ws.Range("G26:J26").ClearContents()

# --- Row 27: "What is the mode of 2/ 5?" -> ": What is the mode of the set {2, 3, 3, 5, 7}?" ---
$ws.Range("A27").Value = ": What is the mode of the set {2, 3, 3, 5, 7}?"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = "B"
$ws.Range("G27:J27").ClearContents()

# --- Column A width ---
$ws.Columns.Item(1).ColumnWidth = 44.17

# --- outlineLevelCol artifact (set + remove outline on a far column so the
#     sheet-level high-water-mark attribute sticks without leaving a visible
#     column definition in the used range) ---
$ws.Columns("Z:Z").OutlineLevel = 5
$ws.Columns("Z:Z").Delete()

# --- Selection ---
$ws.Range("F27").Select()
